# Add debug info for check_elimination func: populate the Week3 (column E)
# results for the players that have already been processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Death Pool")

$ws.Range("E2").Value = "49ers"
$ws.Range("E3").Value = "Dolphins"
$ws.Range("E5").Value = "Packers"
$ws.Range("E7").Value = "Ravens"
$ws.Range("E8").Value = "Ravens"
$ws.Range("E9").Value = "Chiefs"
$ws.Range("E10").Value = "Jaguars"

# Move the active selection to match where the author left off editing.
$ws.Range("E11").Select()
